# Add a new "text" column (D) with an example of the text the program will
# write, alongside the existing command/alternates/description/jargon columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D, matching the bold style used by the other headers
# in row 1 (copy the formatting from A1, then overwrite the value).
$null = $ws.Range("A1").Copy()
$null = $ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "text"

# Example data row for the new column.
$ws.Range("D2").Value = "this is what the program will write"

# Size the new column similarly to the other descriptive columns (~33.29
# characters wide, matching the width Excel stored for column D).
$ws.Columns.Item(4).ColumnWidth = 32.5

# Move the selection to the newly added column, as a user would after typing
# the new values in.
$null = $ws.Range("D5").Select()
